# Refresh Market Board snapshot columns (H:N) across the per-job leve-profit
# tables. Scheduled runner pulls fresh Universalis prices; this just pokes the
# updated currentAveragePrice* / LevePrice* / LeveProfit* cells per row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 1232.1111
$ws.Range("I19").Value = 1150.3334
$ws.Range("K19").Value = 1150.3334
$ws.Range("M19").Value = -975.3334

# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 500
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1152

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 10907.571
$ws.Range("J86").Value = 6720.8
$ws.Range("L86").Value = 6720.8
$ws.Range("N86").Value = -8966.799999999999

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 10907.571
$ws.Range("J89").Value = 6720.8
$ws.Range("L89").Value = 33604
$ws.Range("N89").Value = -44836

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 1535.619
$ws.Range("I98").Value = 1112.4
$ws.Range("K98").Value = 1112.4
$ws.Range("M98").Value = 385.5999999999999

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 2036.25
$ws.Range("I106").Value = 2127.1428
$ws.Range("K106").Value = 2127.1428
$ws.Range("M106").Value = -1496.1428

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 1535.619
$ws.Range("I122").Value = 1112.4
$ws.Range("K122").Value = 3337.2
$ws.Range("M122").Value = -887.2000000000003

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 3598.7
$ws.Range("I132").Value = 3623.5
$ws.Range("K132").Value = 10870.5
$ws.Range("M132").Value = -8340.5

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate / Bronze Plate
$ws.Range("H4").Value = 299.66666
$ws.Range("I4").Value = 199
$ws.Range("K4").Value = 199
$ws.Range("M4").Value = -83

# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 2211.25
$ws.Range("I5").Value = 1058.2
$ws.Range("K5").Value = 1058.2
$ws.Range("M5").Value = -946.2

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 4793.788
$ws.Range("I32").Value = 3167.6128
$ws.Range("K32").Value = 3167.6128
$ws.Range("M32").Value = -2880.6128

# Row 39: Aurochs Star / Bull Hoplon
$ws.Range("H39").Value = 19999.5
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4480

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 5843
$ws.Range("I45").Value = 6666.3076
$ws.Range("K45").Value = 6666.3076
$ws.Range("M45").Value = -6289.3076

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 31255656
$ws.Range("I74").Value = 34488404
$ws.Range("J74").Value = 5766.3335
$ws.Range("K74").Value = 34488404
$ws.Range("L74").Value = 5766.3335
$ws.Range("M74").Value = -34487530
$ws.Range("N74").Value = -7514.3335

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 31255656
$ws.Range("I77").Value = 34488404
$ws.Range("J77").Value = 5766.3335
$ws.Range("K77").Value = 172442020
$ws.Range("L77").Value = 28831.6675
$ws.Range("M77").Value = -172437652
$ws.Range("N77").Value = -37567.6675

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 150000
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 2211.25
$ws.Range("I4").Value = 1058.2
$ws.Range("K4").Value = 1058.2
$ws.Range("M4").Value = -943.2

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1964.9412
$ws.Range("I86").Value = 2199.7856
$ws.Range("J86").Value = 869
$ws.Range("K86").Value = 2199.7856
$ws.Range("L86").Value = 869
$ws.Range("M86").Value = -1076.7856
$ws.Range("N86").Value = -3115

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1964.9412
$ws.Range("I89").Value = 2199.7856
$ws.Range("J89").Value = 869
$ws.Range("K89").Value = 10998.928
$ws.Range("L89").Value = 4345
$ws.Range("M89").Value = -5382.928
$ws.Range("N89").Value = -15577

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1736.1
$ws.Range("I99").Value = 1720.625
$ws.Range("K99").Value = 1720.625
$ws.Range("M99").Value = -222.625

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 22733068
$ws.Range("I134").Value = 22733068
$ws.Range("K134").Value = 68199204
$ws.Range("M134").Value = -68196669

$ws = $wb.Worksheets.Item("CRP")
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 100003384
$ws.Range("J132").Value = 3554
$ws.Range("L132").Value = 10662
$ws.Range("N132").Value = -15722

$ws = $wb.Worksheets.Item("CUL")
# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1900
$ws.Range("J132").Value = 1900
$ws.Range("L132").Value = 17100
$ws.Range("N132").Value = -22160

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 930.3
$ws.Range("I134").Value = 930.3
$ws.Range("K134").Value = 2790.9
$ws.Range("M134").Value = 2279.1

# Row 136: Simple Is Hardest / Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value = 1413.3334
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 13: A Needle Is a Small Sword / Copper Needle
$ws.Range("H13").Value = 149.75
$ws.Range("J13").Value = 149.75
$ws.Range("L13").Value = 149.75
$ws.Range("N13").Value = -427.75

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3826.9092
$ws.Range("I126").Value = 2566.5557
$ws.Range("J126").Value = 9498.5
$ws.Range("K126").Value = 7699.6671
$ws.Range("L126").Value = 28495.5
$ws.Range("M126").Value = -5229.6671
$ws.Range("N126").Value = -33435.5

$ws = $wb.Worksheets.Item("LTW")
# Row 25: A Rush on Ringbands / Hard Leather Ringbands
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 5349
$ws.Range("I61").Value = 5421.4443
$ws.Range("J61").Value = 4697
$ws.Range("K61").Value = 5421.4443
$ws.Range("L61").Value = 4697
$ws.Range("M61").Value = -5219.4443
$ws.Range("N61").Value = -5101

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 5349
$ws.Range("I113").Value = 5421.4443
$ws.Range("J113").Value = 4697
$ws.Range("K113").Value = 5421.4443
$ws.Range("L113").Value = 4697
$ws.Range("M113").Value = -3251.4443
$ws.Range("N113").Value = -9037

$ws = $wb.Worksheets.Item("WVR")
# Row 29: Getting Handsy / Cotton Dress Gloves
$ws.Range("H29").Value = 11199.6
$ws.Range("J29").Value = 9999.333000000001
$ws.Range("L29").Value = 9999.333000000001
$ws.Range("N29").Value = -10579.333

# Row 40: Our Man in Ul'dah / Velveteen Work Gloves
$ws.Range("H40").Value = 19000
$ws.Range("I40").Value = 19000
$ws.Range("K40").Value = 19000
$ws.Range("M40").Value = -18851

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 20009140
$ws.Range("I132").Value = 31253296
$ws.Range("J132").Value = 19531.223
$ws.Range("K132").Value = 93759888
$ws.Range("L132").Value = 58593.66900000001
$ws.Range("M132").Value = -93757358
$ws.Range("N132").Value = -63653.66900000001
